$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab/name
$ws.Name = "Through 2021-11-06"

# Update row 13 label
$ws.Range("A13").Value = "November (through 11-06)"

# Row 13 updates
$ws.Range("C13").Value = 7
$ws.Range("F13").Value = 15
$ws.Range("I13").Value = 20
$ws.Range("J13").Value = 0.0476
$ws.Range("K13").Value = 4
$ws.Range("M13").Value = 0.2222
$ws.Range("O13").Value = 9
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = 42
$ws.Range("S13").Value = 0.0233
$ws.Range("U13").Value = 42

# Row 14 (Total) updates
$ws.Range("C14").Value = 233
$ws.Range("D14").Value = 0.1208
$ws.Range("F14").Value = 449
$ws.Range("G14").Value = 0.1038
$ws.Range("I14").Value = 669
$ws.Range("J14").Value = 0.0848
$ws.Range("K14").Value = 70
$ws.Range("M14").Value = 0.1106
$ws.Range("O14").Value = 443
$ws.Range("P14").Value = 0.0978
$ws.Range("Q14").Value = 55
$ws.Range("R14").Value = 1045
$ws.Range("S14").Value = 0.05
$ws.Range("U14").Value = 1403
$ws.Range("V14").Value = 0.0559
